$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (L3) and formulas (L4, then L5:L8 as one shared-formula fill)
$ws.Range("L3").Value = "1/wzm RMS"
$ws.Range("L4").Formula = "=1/I4"
$ws.Range("L5:L8").Formula = "=1/I5"

# New data table rows 12-16 (A:B) - pasted values from the chart source data
$ws.Range("A12").Value = 0.3
$ws.Range("B12").Value = 0.030204778156996601

$ws.Range("A13").Value = 3
$ws.Range("B13").Value = 0.29969418960244598

$ws.Range("A14").Value = 1000
$ws.Range("B14").Value = 89.102564102564102

$ws.Range("A15").Value = 30
$ws.Range("B15").Value = 2.6347826086956498

$ws.Range("A16").Value = 300
$ws.Range("B16").Value = 27.951807228915701

$ws.Range("B13").Select()
